$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '66.582.83'
$ws.Range('E2').Value = '  -4.70%  '
$ws.Range('D3').Value = '3.462.99'
$ws.Range('E3').Value = '  -6.17%  '
Set-TextValue 'D4' '1.00'
$ws.Range('E4').Value = '  +0.18%  '
Set-TextValue 'D5' '602.69'
$ws.Range('E5').Value = '  -7.03%  '
Set-TextValue 'D6' '148.31'
$ws.Range('E6').Value = '  -8.24%  '
$ws.Range('D7').Value = '3.460.41'
$ws.Range('E7').Value = '  -6.19%  '
$ws.Range('E8').Value = '  +0.05%  '
Set-TextValue 'D9' '0.478'
$ws.Range('E9').Value = '  -4.91%  '
Set-TextValue 'D10' '0.137'
$ws.Range('E10').Value = '  -6.05%  '
Set-TextValue 'D11' '6.88'
$ws.Range('E11').Value = '  -4.21%  '
Set-TextValue 'D12' '0.421'
$ws.Range('E12').Value = '  -5.42%  '
Set-TextValue 'D13' '0.0000217'
$ws.Range('E13').Value = '  -7.10%  '
$ws.Range('D14').Value = '4.056.86'
$ws.Range('E14').Value = '  -5.88%  '
Set-TextValue 'D15' '31.27'
$ws.Range('E15').Value = '  -4.76%  '
$ws.Range('D16').Value = '3.471.51'
$ws.Range('E16').Value = '  -5.84%  '
$ws.Range('D17').Value = '66.596.29'
$ws.Range('E17').Value = '  -4.63%  '
$ws.Range('E18').Value = '  -0.54%  '
Set-TextValue 'D19' '6.35'
$ws.Range('E19').Value = '  -2.90%  '
Set-TextValue 'D20' '14.97'
$ws.Range('E20').Value = '  -6.81%  '
Set-TextValue 'D21' '441.28'
$ws.Range('E21').Value = '  -6.51%  '
Set-TextValue 'D22' '8.97'
$ws.Range('E22').Value = '  -13.96%  '
Set-TextValue 'D23' '0.621'
$ws.Range('E23').Value = '  -4.92%  '
Set-TextValue 'D24' '77.05'
$ws.Range('E24').Value = '  -3.87%  '
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('D26').Value = '3.603.50'
$ws.Range('E26').Value = '  -6.06%  '
Set-TextValue 'D27' '0.0000123'
$ws.Range('E27').Value = '  -4.20%  '
Set-TextValue 'D28' '10.01'
$ws.Range('E28').Value = '  -8.97%  '
Set-TextValue 'D29' '8.21'
$ws.Range('E29').Value = '  -10.59%  '
$ws.Range('E30').Value = '  -5.77%  '
$ws.Range('E31').Value = '  -8.98%  '
Set-TextValue 'D32' '0.999'
$ws.Range('E32').Value = '  -0.13%  '
Set-TextValue 'D33' '0.159'
Set-TextValue 'D34' '25.44'
$ws.Range('E34').Value = '  -5.23%  '
Set-TextValue 'D35' '6.11'
$ws.Range('E35').Value = '  -6.81%  '
$ws.Range('E36').Value = '  -8.43%  '
$ws.Range('D37').Value = '3.453.10'
$ws.Range('E37').Value = '  -6.37%  '
Set-TextValue 'D38' '7.89'
$ws.Range('E38').Value = '  -6.69%  '
$ws.Range('E39').Value = '  +0.04%  '
Set-TextValue 'D40' '1.00'
$ws.Range('E40').Value = '  +0.24%  '
Set-TextValue 'D41' '173.21'
$ws.Range('E41').Value = '  -3.11%  '
Set-TextValue 'D42' '2.16'
Set-TextValue 'D43' '0.0860'
$ws.Range('E43').Value = '  -5.02%  '
Set-TextValue 'D44' '5.45'
$ws.Range('E44').Value = '  -7.72%  '
Set-TextValue 'D45' '0.877'
$ws.Range('E45').Value = '  -5.96%  '
Set-TextValue 'D46' '45.21'
$ws.Range('E46').Value = '  -3.63%  '
Set-TextValue 'D47' '1.22'
$ws.Range('E47').Value = '  -3.81%  '
Set-TextValue 'D48' '26.00'
$ws.Range('E48').Value = '  -11.31%  '
Set-TextValue 'D49' '2.48'
$ws.Range('E49').Value = '  -13.41%  '
Set-TextValue 'D50' '7.52'
$ws.Range('E50').Value = '  -4.44%  '
Set-TextValue 'D51' '1.00'
$ws.Range('E51').Value = '  -5.15%  '
